$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.024701
$ws.Range("N2").Value = 9.074103000000001
$ws.Range("O2").Value = 0.1596375877334842
$ws.Range("P2").Value = 0.1596375877334843
$ws.Range("Q2").Value = 31.56046624457667
$ws.Range("R2").Value = 284.04419620119
$ws.Range("S2").Value = 0.1550255297752276
$ws.Range("T2").Value = 0.1550255297752277

# Row 3
$ws.Range("O3").Value = 0.6072559333217162
$ws.Range("P3").Value = 0.6072559333217163
$ws.Range("S3").Value = 0.5897118223154115
$ws.Range("T3").Value = 0.5897118223154116

# Row 4
$ws.Range("M4").Value = 4.368554666666666
$ws.Range("N4").Value = 13.105664
$ws.Range("O4").Value = 0.2305634602787257
$ws.Range("P4").Value = 0.2305634602787257
$ws.Range("Q4").Value = 45.58256240696888
$ws.Range("R4").Value = 410.2430616627199
$ws.Range("S4").Value = 0.2239022969715165
$ws.Range("T4").Value = 0.2239022969715165

# Row 5
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.04818333333333333
$ws.Range("N5").Value = 0.14455
$ws.Range("O5").Value = 0.002543018666073676
$ws.Range("P5").Value = 0.002543018666073677
$ws.Range("Q5").Value = 0.5027566246111111
$ws.Range("R5").Value = 4.5248096215
$ws.Range("S5").Value = 0.002469548817002535
$ws.Range("T5").Value = 0.002469548817002536

# Row 6
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.024701
$ws.Range("N6").Value = 9.074103000000001
$ws.Range("O6").Value = 0.1596375877334842
$ws.Range("P6").Value = 0.1596375877334843
$ws.Range("Q6").Value = 0.9389337338220001
$ws.Range("R6").Value = 8.450403604398
$ws.Range("S6").Value = 0.004612057958256584
$ws.Range("T6").Value = 0.004612057958256586

# Row 7
$ws.Range("O7").Value = 0.6072559333217162
$ws.Range("P7").Value = 0.6072559333217163
$ws.Range("S7").Value = 0.01754411100630469
$ws.Range("T7").Value = 0.01754411100630469

# Row 8
$ws.Range("M8").Value = 4.368554666666666
$ws.Range("N8").Value = 13.105664
$ws.Range("O8").Value = 0.2305634602787257
$ws.Range("P8").Value = 0.2305634602787257
$ws.Range("Q8").Value = 1.356095476736
$ws.Range("R8").Value = 12.204859290624
$ws.Range("S8").Value = 0.006661163307209188
$ws.Range("T8").Value = 0.00666116330720919

# Row 9
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.04818333333333333
$ws.Range("N9").Value = 0.14455
$ws.Range("O9").Value = 0.002543018666073676
$ws.Range("P9").Value = 0.002543018666073677
$ws.Range("Q9").Value = 0.0149571667
$ws.Range("R9").Value = 0.1346145003
$ws.Range("S9").Value = 0.00007346984907114117
$ws.Range("T9").Value = 0.00007346984907114119
